$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the original first data row (old row 2); everything below shifts up.
$ws.Rows.Item(2).Delete()

# Append the new accelerometer samples captured on May 9th (rows 21-31).
$ws.Cells.Item(21, 1).Value = -5.779280513524995
$ws.Cells.Item(21, 2).Value = -3.271673738956447
$ws.Cells.Item(21, 3).Value = 4.346601516008363
$ws.Cells.Item(22, 1).Value = 3.140387788414934
$ws.Cells.Item(22, 2).Value = 0.259726375341407
$ws.Cells.Item(22, 3).Value = -4.785581156611421
$ws.Cells.Item(23, 1).Value = -1.632258296012878
$ws.Cells.Item(23, 2).Value = 0.6425724923610687
$ws.Cells.Item(23, 3).Value = -3.22618693113327
$ws.Cells.Item(24, 1).Value = -3.271841421723368
$ws.Cells.Item(24, 2).Value = 0.07577018067240526
$ws.Cells.Item(24, 3).Value = -1.077775649726385
$ws.Cells.Item(25, 1).Value = -3.740465611219407
$ws.Cells.Item(25, 2).Value = -0.2502757757902145
$ws.Cells.Item(25, 3).Value = -5.230584308505059
$ws.Cells.Item(26, 1).Value = -1.173786669969556
$ws.Cells.Item(26, 2).Value = -1.206141140311958
$ws.Cells.Item(26, 3).Value = -5.999948702752588
$ws.Cells.Item(27, 1).Value = -0.527452439069747
$ws.Cells.Item(27, 2).Value = -1.933494433760643
$ws.Cells.Item(27, 3).Value = -1.995455801486972
$ws.Cells.Item(28, 1).Value = -2.818732134997842
$ws.Cells.Item(28, 2).Value = -1.367858927696945
$ws.Cells.Item(28, 3).Value = 0.861171409487724
$ws.Cells.Item(29, 1).Value = -4.338251754641532
$ws.Cells.Item(29, 2).Value = -0.358771674335002
$ws.Cells.Item(29, 3).Value = 0.811524987220763
$ws.Cells.Item(30, 1).Value = -5.003720842301852
$ws.Cells.Item(30, 2).Value = -0.2583636995404971
$ws.Cells.Item(30, 3).Value = 1.431181490421301
$ws.Cells.Item(31, 1).Value = -6.261017680168152
$ws.Cells.Item(31, 2).Value = -1.120874315500259
$ws.Cells.Item(31, 3).Value = 4.274679899215698
